$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 527873
$ws.Range("E2").Value = -32740
$ws.Range("F2").Value = -32495
$ws.Range("G2").Value = -31155
$ws.Range("H2").Value = -22138
$ws.Range("I2").Value = -17688
$ws.Range("J2").Value = -4450
$ws.Range("K2").Value = 531138
$ws.Range("L2").Value = 365355
$ws.Range("M2").Value = 165782
$ws.Range("N2").Value = 151742
$ws.Range("O2").Value = 14041
$ws.Range("P2").Value = 3800
$ws.Range("Q2").Value = 15024
$ws.Range("R2").Value = -4281
$ws.Range("S2").Value = 8368
$ws.Range("T2").Value = 12997
$ws.Range("U2").Value = 2026
$ws.Range("V2").Value = 190164
$ws.Range("W2").Value = -6.2
$ws.Range("X2").Value = -4.19
$ws.Range("Y2").Value = -10.91
$ws.Range("Z2").Value = -4.16
$ws.Range("AA2").Value = 220.38
$ws.Range("AB2").Value = 4047.8
$ws.Range("AC2").Value = -21964
$ws.Range("AD2").Value = -4.42
$ws.Range("AE2").Value = 230533
$ws.Range("AF2").Value = 0.42
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 80533716

# Row 3
$ws.Range("D3").Value = 463176
$ws.Range("E3").Value = -15849
$ws.Range("F3").Value = -15401
$ws.Range("G3").Value = -18600
$ws.Range("H3").Value = -13758
$ws.Range("I3").Value = -13499
$ws.Range("J3").Value = -259
$ws.Range("K3").Value = 494691
$ws.Range("L3").Value = 340430
$ws.Range("M3").Value = 154261
$ws.Range("N3").Value = 137366
$ws.Range("O3").Value = 16895
$ws.Range("P3").Value = 3800
$ws.Range("Q3").Value = -6374
$ws.Range("R3").Value = -1386
$ws.Range("S3").Value = 6127
$ws.Range("T3").Value = 12724
$ws.Range("U3").Value = -19099
$ws.Range("V3").Value = 182071
$ws.Range("W3").Value = -3.42
$ws.Range("X3").Value = -2.97
$ws.Range("Y3").Value = -9.34
$ws.Range("Z3").Value = -2.68
$ws.Range("AA3").Value = 220.68
$ws.Range("AB3").Value = 3669.5
$ws.Range("AC3").Value = -16761
$ws.Range("AD3").Value = -4.43
$ws.Range("AE3").Value = 195188
$ws.Range("AF3").Value = 0.38
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 80533716

# Row 4
$ws.Range("D4").Value = 223004
$ws.Range("E4").Value = 3915
$ws.Range("F4").Value = 16419
$ws.Range("G4").Value = 960
$ws.Range("H4").Value = 6270
$ws.Range("I4").Value = 5452
$ws.Range("J4").Value = 818
$ws.Range("K4").Value = 489493
$ws.Range("L4").Value = 311627
$ws.Range("M4").Value = 177866
$ws.Range("N4").Value = 159369
$ws.Range("O4").Value = 18496
$ws.Range("P4").Value = 3800
$ws.Range("Q4").Value = 26492
$ws.Range("R4").Value = -5248
$ws.Range("S4").Value = -8759
$ws.Range("T4").Value = 11355
$ws.Range("U4").Value = 15138
$ws.Range("V4").Value = 179808
$ws.Range("W4").Value = 1.76
$ws.Range("X4").Value = 2.81
$ws.Range("Y4").Value = 3.67
$ws.Range("Z4").Value = 1.27
$ws.Range("AA4").Value = 175.2
$ws.Range("AB4").Value = 3831.27
$ws.Range("AC4").Value = 6769
$ws.Range("AD4").Value = 18.16
$ws.Range("AE4").Value = 226454
$ws.Range("AF4").Value = 0.54
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 80533716

# Row 5
$ws.Range("D5").Value = 154688
$ws.Range("E5").Value = 146
$ws.Range("F5").Value = 146
$ws.Range("G5").Value = -269
$ws.Range("H5").Value = 26931
$ws.Range("I5").Value = 24578
$ws.Range("J5").Value = 2354
$ws.Range("K5").Value = 304088
$ws.Range("L5").Value = 180368
$ws.Range("M5").Value = 123720
$ws.Range("N5").Value = 111214
$ws.Range("O5").Value = 12506
$ws.Range("P5").Value = 2833
$ws.Range("Q5").Value = 5933
$ws.Range("R5").Value = 5963
$ws.Range("S5").Value = -21225
$ws.Range("T5").Value = 3733
$ws.Range("U5").Value = 2200
$ws.Range("V5").Value = 53366
$ws.Range("W5").Value = 0.1
$ws.Range("X5").Value = 17.41
$ws.Range("Y5").Value = 18.17
$ws.Range("Z5").Value = 6.79
$ws.Range("AA5").Value = 145.79
$ws.Range("AB5").Value = 6101.36
$ws.Range("AC5").Value = 37625
$ws.Range("AD5").Value = 2.55
$ws.Range("AE5").Value = 185395
$ws.Range("AF5").Value = 0.52
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 60045754

# Row 6
$ws.Range("D6").Value = 131199
$ws.Range("E6").Value = -5225
$ws.Range("F6").Value = -5225
$ws.Range("G6").Value = -8540
$ws.Range("H6").Value = -4536
$ws.Range("I6").Value = -4891
$ws.Range("K6").Value = 247299
$ws.Range("L6").Value = 116189
$ws.Range("M6").Value = 131110
$ws.Range("N6").Value = 118822
$ws.Range("P6").Value = 3539
$ws.Range("Q6").Value = 1265
$ws.Range("R6").Value = -57
$ws.Range("S6").Value = -3294
$ws.Range("T6").Value = 2258
$ws.Range("U6").Value = -992
$ws.Range("V6").Value = 39004
$ws.Range("W6").Value = -3.98
$ws.Range("X6").Value = -3.46
$ws.Range("Y6").Value = -4.25
$ws.Range("Z6").Value = -1.65
$ws.Range("AA6").Value = 88.62
$ws.Range("AB6").Value = 5235.63
$ws.Range("AC6").Value = -7296
$ws.Range("AD6").Value = -17.61
$ws.Range("AE6").Value = 168031
$ws.Range("AF6").Value = 0.76
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 70773116

# Row 7
$ws.Range("D7").Value = 148934
$ws.Range("E7").Value = 1613
$ws.Range("G7").Value = 1293
$ws.Range("H7").Value = 2509
$ws.Range("I7").Value = 1935
$ws.Range("K7").Value = 253992
$ws.Range("L7").Value = 121219
$ws.Range("M7").Value = 132774
$ws.Range("N7").Value = 118390
$ws.Range("P7").Value = 3540
$ws.Range("Q7").Value = 7280
$ws.Range("R7").Value = 2771
$ws.Range("S7").Value = -6592
$ws.Range("T7").Value = 3316
$ws.Range("U7").Value = 1957
$ws.Range("W7").Value = 1.08
$ws.Range("X7").Value = 1.69
$ws.Range("Y7").Value = 1.63
$ws.Range("Z7").Value = 1
$ws.Range("AA7").Value = 91.3
$ws.Range("AC7").Value = 2735
$ws.Range("AD7").Value = 41.69
$ws.Range("AE7").Value = 167419
$ws.Range("AF7").Value = 0.68
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 159719
$ws.Range("E8").Value = 2955
$ws.Range("G8").Value = 2966
$ws.Range("H8").Value = 2205
$ws.Range("I8").Value = 1750
$ws.Range("K8").Value = 259855
$ws.Range("L8").Value = 125320
$ws.Range("M8").Value = 134536
$ws.Range("N8").Value = 119687
$ws.Range("P8").Value = 3540
$ws.Range("Q8").Value = 6761
$ws.Range("R8").Value = 211
$ws.Range("S8").Value = -7290
$ws.Range("T8").Value = 2994
$ws.Range("U8").Value = 3246
$ws.Range("W8").Value = 1.85
$ws.Range("X8").Value = 1.38
$ws.Range("Y8").Value = 1.47
$ws.Range("Z8").Value = 0.86
$ws.Range("AA8").Value = 93.15000000000001
$ws.Range("AC8").Value = 2472
$ws.Range("AD8").Value = 46.11
$ws.Range("AE8").Value = 169254
$ws.Range("AF8").Value = 0.67
$ws.Range("AG8").Value = 50
$ws.Range("AH8").Value = 0.04
$ws.Range("AI8").Value = 2.02

# Row 9
$ws.Range("D9").Value = 170806
$ws.Range("E9").Value = 4670
$ws.Range("G9").Value = 4217
$ws.Range("H9").Value = 3093
$ws.Range("I9").Value = 2501
$ws.Range("K9").Value = 265109
$ws.Range("L9").Value = 128315
$ws.Range("M9").Value = 136793
$ws.Range("N9").Value = 121662
$ws.Range("P9").Value = 3540
$ws.Range("Q9").Value = 5626
$ws.Range("R9").Value = 782
$ws.Range("S9").Value = -6928
$ws.Range("T9").Value = 3090
$ws.Range("U9").Value = 3538
$ws.Range("W9").Value = 2.73
$ws.Range("X9").Value = 1.81
$ws.Range("Y9").Value = 2.07
$ws.Range("Z9").Value = 1.18
$ws.Range("AA9").Value = 93.8
$ws.Range("AC9").Value = 3533
$ws.Range("AD9").Value = 32.26
$ws.Range("AE9").Value = 172047
$ws.Range("AF9").Value = 0.66
$ws.Range("AG9").Value = 111
$ws.Range("AH9").Value = 0.1
$ws.Range("AI9").Value = 3.15
